# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-job Leve profit sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0
$ws.Range("H40").Value = 3679.0527
$ws.Range("J40").Value = 3933.6667
$ws.Range("L40").Value = 3933.6667
$ws.Range("N40").Value = -4283.6667

# Hunk 1
$ws.Range("H51").Value = 9332.5
$ws.Range("J51").Value = 7498.5
$ws.Range("L51").Value = 7498.5
$ws.Range("N51").Value = -8466.5

# Hunk 2
$ws.Range("H99").Value = 2036.3334
$ws.Range("J99").Value = 4556.4
$ws.Range("L99").Value = 13669.2
$ws.Range("N99").Value = -16665.2

# Hunk 3
$ws.Range("H111").Value = 975.3333
$ws.Range("I111").Value = 975.3333
$ws.Range("K111").Value = 2925.9999
$ws.Range("M111").Value = 141.0001000000002

# Hunk 4
$ws.Range("H118").Value = 1003.8333
$ws.Range("I118").Value = 506.8
$ws.Range("K118").Value = 1520.4
$ws.Range("M118").Value = 136.5999999999999

# Hunk 5
$ws.Range("H137").Value = 5287.533
$ws.Range("I137").Value = 1274.5714
$ws.Range("K137").Value = 3823.7142
$ws.Range("M137").Value = -1273.7142

# Hunk 6
$ws.Range("H141").Value = 4195.1875
$ws.Range("I141").Value = 4480.2144
$ws.Range("K141").Value = 13440.6432
$ws.Range("M141").Value = -8260.643199999999

$ws = $wb.Worksheets.Item("ARM")
# Hunk 7
$ws.Range("H32").Value = 44889.074
$ws.Range("I32").Value = 51530
$ws.Range("J32").Value = 15669
$ws.Range("K32").Value = 51530
$ws.Range("L32").Value = 15669
$ws.Range("M32").Value = -51243
$ws.Range("N32").Value = -16243

# Hunk 8
$ws.Range("H63").Value = 4177.6
$ws.Range("I63").Value = 2629.3333
$ws.Range("J63").Value = 6500
$ws.Range("K63").Value = 2629.3333
$ws.Range("L63").Value = 6500
$ws.Range("M63").Value = -1943.3333
$ws.Range("N63").Value = -7872

# Hunk 9
$ws.Range("H66").Value = 4177.6
$ws.Range("I66").Value = 2629.3333
$ws.Range("J66").Value = 6500
$ws.Range("K66").Value = 13146.6665
$ws.Range("L66").Value = 32500
$ws.Range("M66").Value = -9714.666499999999
$ws.Range("N66").Value = -39364

# Hunk 10
$ws.Range("H110").Value = 2783
$ws.Range("J110").Value = 2119.6
$ws.Range("L110").Value = 2119.6
$ws.Range("N110").Value = -6209.6

# Hunk 11
$ws.Range("H132").Value = 62587.234
$ws.Range("I132").Value = 75236.36
$ws.Range("J132").Value = 3558
$ws.Range("K132").Value = 225709.08
$ws.Range("L132").Value = 10674
$ws.Range("M132").Value = -223179.08
$ws.Range("N132").Value = -15734

$ws = $wb.Worksheets.Item("BSM")
# Hunk 12
$ws.Range("H86").Value = 2283.0625
$ws.Range("I86").Value = 2328.1428
$ws.Range("J86").Value = 2248
$ws.Range("K86").Value = 2328.1428
$ws.Range("L86").Value = 2248
$ws.Range("M86").Value = -1205.1428
$ws.Range("N86").Value = -4494

# Hunk 13
$ws.Range("H89").Value = 2283.0625
$ws.Range("I89").Value = 2328.1428
$ws.Range("J89").Value = 2248
$ws.Range("K89").Value = 11640.714
$ws.Range("L89").Value = 11240
$ws.Range("M89").Value = -6024.714
$ws.Range("N89").Value = -22472

# Hunk 14
$ws.Range("H99").Value = 80915.234
$ws.Range("I99").Value = 114655.445
$ws.Range("K99").Value = 114655.445
$ws.Range("M99").Value = -113157.445

# Hunk 15
$ws.Range("H134").Value = 2443.5
$ws.Range("I134").Value = 2443.5
$ws.Range("K134").Value = 7330.5
$ws.Range("M134").Value = -4795.5

$ws = $wb.Worksheets.Item("CRP")
# Hunk 16
$ws.Range("H47").Value = 37999.668
$ws.Range("I47").Value = 34000
$ws.Range("J47").Value = 39999.5
$ws.Range("K47").Value = 34000
$ws.Range("L47").Value = 39999.5
$ws.Range("M47").Value = -33434
$ws.Range("N47").Value = -41131.5

# Hunk 17
$ws.Range("H58").Value = 103331.8
$ws.Range("I58").Value = 103331.8
$ws.Range("K58").Value = 103331.8
$ws.Range("M58").Value = -103128.8

# Hunk 18
$ws.Range("H136").Value = 103331.8
$ws.Range("I136").Value = 103331.8
$ws.Range("K136").Value = 309995.4
$ws.Range("M136").Value = -307445.4

$ws = $wb.Worksheets.Item("CUL")
# Hunk 19
$ws.Range("H96").Value = 20007.666
$ws.Range("J96").Value = 24999
$ws.Range("L96").Value = 74997
$ws.Range("N96").Value = -79115

# Hunk 20
$ws.Range("H99").Value = 607.5
$ws.Range("I99").Value = 607.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1822.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 423.5
$ws.Range("N99").ClearContents()

# Hunk 21
$ws.Range("H128").Value = 154596.8
$ws.Range("I128").Value = 154596.8
$ws.Range("K128").Value = 463790.4
$ws.Range("M128").Value = -458810.4

$ws = $wb.Worksheets.Item("GSM")
# Hunk 22
$ws.Range("H2").Value = 66.76470999999999
$ws.Range("I2").Value = 64.416664
$ws.Range("K2").Value = 64.416664
$ws.Range("M2").Value = 48.583336

# Hunk 23
$ws.Range("H36").Value = 5000
$ws.Range("J36").Value = 5000
$ws.Range("L36").Value = 5000
$ws.Range("N36").Value = -5970

# Hunk 24
$ws.Range("H80").Value = 3158.875
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002

# Hunk 25
$ws.Range("H83").Value = 3158.875
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008

# Hunk 26
$ws.Range("H111").Value = 28000
$ws.Range("J111").Value = 28000
$ws.Range("L111").Value = 28000
$ws.Range("N111").Value = -34134

# Hunk 27
$ws.Range("H126").Value = 19999
$ws.Range("I126").Value = 14998
$ws.Range("K126").Value = 44994
$ws.Range("M126").Value = -42524

# Hunk 28
$ws.Range("H139").Value = 69004.5
$ws.Range("J139").Value = 69004.5
$ws.Range("L139").Value = 69004.5
$ws.Range("N139").Value = -79284.5

$ws = $wb.Worksheets.Item("LTW")
# Hunk 29
$ws.Range("H22").Value = 2625.2
$ws.Range("I22").Value = 723.5
$ws.Range("J22").Value = 4798.5713
$ws.Range("K22").Value = 723.5
$ws.Range("L22").Value = 4798.5713
$ws.Range("M22").Value = -428.5
$ws.Range("N22").Value = -5388.5713

# Hunk 30
$ws.Range("H27").Value = 2625.2
$ws.Range("I27").Value = 723.5
$ws.Range("J27").Value = 4798.5713
$ws.Range("K27").Value = 723.5
$ws.Range("L27").Value = 4798.5713
$ws.Range("M27").Value = -616.5
$ws.Range("N27").Value = -5012.5713

# Hunk 31
$ws.Range("H40").Value = 4107.1665
$ws.Range("I40").Value = 3235.75
$ws.Range("J40").Value = 5850
$ws.Range("K40").Value = 3235.75
$ws.Range("L40").Value = 5850
$ws.Range("M40").Value = -3099.75
$ws.Range("N40").Value = -6122

# Hunk 32
$ws.Range("H47").Value = 10001
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 10001
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 10001
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -10981

# Hunk 33
$ws.Range("H52").Value = 10001
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 10001
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 10001
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -10467

# Hunk 34
$ws.Range("H93").Value = 1483.238
$ws.Range("I93").Value = 903.3570999999999
$ws.Range("K93").Value = 903.3570999999999
$ws.Range("M93").Value = 344.6429000000001

# Hunk 35
$ws.Range("H100").Value = 3296.6
$ws.Range("I100").Value = 2998.375
$ws.Range("J100").Value = 4489.5
$ws.Range("K100").Value = 2998.375
$ws.Range("L100").Value = 4489.5
$ws.Range("M100").Value = -2457.375
$ws.Range("N100").Value = -5571.5

# Hunk 36
$ws.Range("H132").Value = 39811.78
$ws.Range("I132").Value = 40773.484
$ws.Range("K132").Value = 122320.452
$ws.Range("M132").Value = -119790.452

# Hunk 37
$ws.Range("H136").Value = 4014.4443
$ws.Range("I136").Value = 2791.3635
$ws.Range("K136").Value = 8374.0905
$ws.Range("M136").Value = -5824.0905

$ws = $wb.Worksheets.Item("WVR")
# Hunk 38
$ws.Range("H126").Value = 73900.14999999999
$ws.Range("I126").Value = 88518.13
$ws.Range("J126").Value = 9581
$ws.Range("K126").Value = 265554.39
$ws.Range("L126").Value = 28743
$ws.Range("M126").Value = -263084.39
$ws.Range("N126").Value = -33683

# Hunk 39
$ws.Range("H136").Value = 2095.975
$ws.Range("I136").Value = 1767.5625
$ws.Range("K136").Value = 5302.6875
$ws.Range("M136").Value = -2752.6875
